$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings like "29.240.02" or "1.000" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.240.02'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '1.845.53'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '243.08'
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").Value = '0.6652'
$ws.Range("E6").Value = '  -1.31%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").Value = '  +6.78%  '

$ws.Range("D9").Value = '0.07470'
$ws.Range("E9").Value = '  +0.13%  '

$ws.Range("D10").Value = '0.2966'
$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("D11").Value = '23.38'
$ws.Range("E11").Value = '  +1.97%  '

$ws.Range("D12").Value = '0.07763'
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").Value = '1.851.71'
$ws.Range("E13").Value = '  +0.91%  '

$ws.Range("D14").Value = '5.031'
$ws.Range("E14").Value = '  -0.14%  '

$ws.Range("D15").Value = '0.6765'
$ws.Range("E15").Value = '  -0.76%  '

$ws.Range("D16").Value = '83.45'
$ws.Range("E16").Value = '  -3.72%  '

$ws.Range("D17").Value = '6.203'
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '0.000008732'
$ws.Range("E18").Value = '  +5.09%  '

$ws.Range("D19").Value = '29.195.54'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = '2.097.06'
$ws.Range("E20").Value = '  +1.87%  '

$ws.Range("D21").Value = '227.59'
$ws.Range("E21").Value = '  -0.91%  '

$ws.Range("D22").Value = '12.56'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").Value = '7.231'
$ws.Range("E24").Value = '  -0.31%  '

$ws.Range("D25").Value = '1.000'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").Value = '158.66'
$ws.Range("E26").Value = '  -1.03%  '

$ws.Range("D27").Value = '8.644'
$ws.Range("E27").Value = '  -0.90%  '

$ws.Range("D28").Value = '0.1403'
$ws.Range("E28").Value = '  -0.57%  '

$ws.Range("D29").Value = '18.08'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("D31").Value = '4.144'
$ws.Range("E31").Value = '  -1.34%  '

$ws.Range("D32").Value = '4.058'
$ws.Range("E32").Value = '  -1.09%  '

$ws.Range("D33").Value = '1.194'
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").Value = '0.05380'
$ws.Range("E34").Value = '  +0.68%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7483'
$ws.Range("E35").Value = '  -1.59%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '1.852'
$ws.Range("E36").Value = '  -2.80%  '

$ws.Range("D37").Value = '1.163'
$ws.Range("E37").Value = '  +1.37%  '

$ws.Range("D38").Value = '2.647'
$ws.Range("E38").Value = '  -1.33%  '

$ws.Range("D39").Value = '1.302.45'
$ws.Range("E39").Value = '  -2.68%  '

$ws.Range("E40").Value = '  -0.73%  '

$ws.Range("D41").Value = '2.758'
$ws.Range("E41").Value = '  +0.59%  '

$ws.Range("D42").Value = '6.409'
$ws.Range("E42").Value = '  +7.28%  '

$ws.Range("D43").Value = '0.9096'
$ws.Range("E43").Value = '  -1.94%  '

$ws.Range("B44").Value = 'XinFinNetwork'
$ws.Range("C44").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D44").Value = '0.08327'
$ws.Range("E44").Value = '  +3.68%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").Value = '  -0.20%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '103.71'
$ws.Range("E46").Value = '  -0.11%  '

$ws.Range("D47").Value = '1.994.89'
$ws.Range("E47").Value = '  +1.73%  '

$ws.Range("D48").Value = '65.51'
$ws.Range("E48").Value = '  +2.32%  '

$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("E51").Value = '  -1.05%  '
